$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.255910158157349
$ws.Range("B1").Value = 1.388409018516541
$ws.Range("C1").Value = 2.308276653289795
$ws.Range("D1").Value = 3.950424671173096
$ws.Range("E1").Value = 1.199276804924011
